$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $row, $col, $text)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws 2 4 "66.624.25"
Set-TextValue $ws 2 5 "  +0.28%  "
Set-TextValue $ws 3 4 "3.225.45"
Set-TextValue $ws 3 5 "  +0.87%  "
Set-TextValue $ws 4 4 "0.999"
Set-TextValue $ws 4 5 "  -0.06%  "
Set-TextValue $ws 5 4 "608.95"
Set-TextValue $ws 5 5 "  +2.08%  "
Set-TextValue $ws 6 4 "158.72"
Set-TextValue $ws 6 5 "  +2.69%  "
Set-TextValue $ws 7 5 "  +0.04%  "
Set-TextValue $ws 8 4 "3.224.37"
Set-TextValue $ws 8 5 "  +0.89%  "
Set-TextValue $ws 9 4 "0.550"
Set-TextValue $ws 9 5 "  +0.65%  "
Set-TextValue $ws 10 5 "  +0.24%  "
Set-TextValue $ws 11 4 "5.70"
Set-TextValue $ws 11 5 "  -5.30%  "
Set-TextValue $ws 12 5 "  -3.03%  "
Set-TextValue $ws 13 4 "0.0000270"
Set-TextValue $ws 13 5 "  +0.81%  "
Set-TextValue $ws 14 4 "38.79"
Set-TextValue $ws 14 5 "  -1.01%  "
Set-TextValue $ws 15 4 "3.757.23"
Set-TextValue $ws 16 4 "66.643.64"
Set-TextValue $ws 16 5 "  +0.32%  "
Set-TextValue $ws 17 4 "7.35"
Set-TextValue $ws 17 5 "  -1.64%  "
Set-TextValue $ws 18 4 "3.231.12"
Set-TextValue $ws 18 5 "  +1.09%  "
Set-TextValue $ws 19 5 "  +1.11%  "
Set-TextValue $ws 20 4 "507.01"
Set-TextValue $ws 21 4 "15.14"
Set-TextValue $ws 21 5 "  -1.77%  "
Set-TextValue $ws 22 4 "0.733"
Set-TextValue $ws 22 5 "  -1.17%  "
Set-TextValue $ws 23 4 "7.99"
Set-TextValue $ws 23 5 "  -1.45%  "
Set-TextValue $ws 24 4 "14.59"
Set-TextValue $ws 24 5 "  -3.09%  "
Set-TextValue $ws 25 5 "  -1.05%  "
Set-TextValue $ws 26 5 "  +0.04%  "
Set-TextValue $ws 27 5 "  -0.21%  "
Set-TextValue $ws 28 5 "  -2.24%  "
Set-TextValue $ws 29 4 "2.35"
Set-TextValue $ws 29 5 "  +0.86%  "
Set-TextValue $ws 30 2 "NEARProtocol"
Set-TextValue $ws 30 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws 30 4 "7.01"
Set-TextValue $ws 30 5 "  -2.96%  "
Set-TextValue $ws 31 2 "Hedera"
Set-TextValue $ws 31 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws 31 4 "0.119"
Set-TextValue $ws 31 5 "  +31.82%  "
Set-TextValue $ws 32 4 "2.92"
Set-TextValue $ws 32 5 "  +0.57%  "
Set-TextValue $ws 33 4 "28.11"
Set-TextValue $ws 33 5 "  -0.71%  "
Set-TextValue $ws 34 5 "  +0.17%  "
Set-TextValue $ws 35 5 "  -4.04%  "
Set-TextValue $ws 36 4 "6.46"
Set-TextValue $ws 36 5 "  -1.30%  "
Set-TextValue $ws 37 4 "55.57"
Set-TextValue $ws 37 5 "  +1.27%  "
Set-TextValue $ws 38 4 "502.16"
Set-TextValue $ws 38 5 "  -2.20%  "
Set-TextValue $ws 39 4 "0.0₃0770"
Set-TextValue $ws 39 5 "  +14.37%  "
Set-TextValue $ws 40 4 "3.10"
Set-TextValue $ws 40 5 "  +7.13%  "
Set-TextValue $ws 41 4 "0.131"
Set-TextValue $ws 41 5 "  +6.35%  "
Set-TextValue $ws 42 4 "0.0419"
Set-TextValue $ws 42 5 "  -1.32%  "
Set-TextValue $ws 43 4 "8.68"
Set-TextValue $ws 43 5 "  -2.45%  "
Set-TextValue $ws 44 5 "  -2.20%  "
Set-TextValue $ws 45 5 "  -0.48%  "
Set-TextValue $ws 46 4 "2.894.25"
Set-TextValue $ws 46 5 "  -0.89%  "
Set-TextValue $ws 47 4 "28.13"
Set-TextValue $ws 47 5 "  -2.33%  "
Set-TextValue $ws 48 4 "2.42"
Set-TextValue $ws 48 5 "  +2.90%  "
Set-TextValue $ws 49 5 "  -0.07%  "
Set-TextValue $ws 50 5 "  -1.20%  "
Set-TextValue $ws 51 4 "122.20"
Set-TextValue $ws 51 5 "  -0.61%  "
